$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of portfolio data for 2025-11-01.
# Column A holds the date as plain text (matches existing rows, which are
# inlineStr/text, not real Excel dates) - force text format before writing
# so Excel doesn't auto-convert the "YYYY-MM-DD" string into a date serial,
# then clear the format override so no extra cell style is introduced.
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "2025-11-01"
$ws.Range("A78").ClearFormats()

$ws.Range("B78").Value = 59.29999923706055
$ws.Range("C78").Value = 410
$ws.Range("D78").Value = 317.75
